$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.953.15'
$ws.Range("E2").Value = '  +2.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.327.73'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.81'
$ws.Range("E5").Value = '  +5.90%  '
$ws.Range("E6").Value = '  +2.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.357.95'
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.39'
$ws.Range("E12").Value = '  +2.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +4.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.752.67'
$ws.Range("E14").Value = '  +0.25%  '
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.820.71'
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.360.83'
$ws.Range("E18").Value = '  +1.30%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.55'
$ws.Range("E19").Value = '  +1.06%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '334.72'
$ws.Range("E20").Value = '  +3.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.20'
$ws.Range("E21").Value = '  +1.84%  '
$ws.Range("E22").Value = '  +2.05%  '
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.76'
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("E25").Value = '  +4.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.48'
$ws.Range("E26").Value = '  -2.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.993'
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("E28").Value = '  +9.75%  '
$ws.Range("E29").Value = '  +4.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.61'
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0737'
$ws.Range("E31").Value = '  +2.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.17'
$ws.Range("E32").Value = '  +1.31%  '
$ws.Range("E33").Value = '  +1.53%  '
$ws.Range("E34").Value = '  +15.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.996'
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.13'
$ws.Range("E38").Value = '  +5.11%  '
$ws.Range("E39").Value = '  +4.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '39.35'
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '149.09'
$ws.Range("E41").Value = '  -1.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.379'
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '282.02'
$ws.Range("E44").Value = '  +1.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.34'
$ws.Range("E45").Value = '  +7.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0505'
$ws.Range("E47").Value = '  +2.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.561'
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("E49").Value = '  +1.94%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.54'
$ws.Range("E50").Value = '  +2.38%  '
$ws.Range("B51").Value = 'Polygon'
$ws.Range("C51").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.382'
$ws.Range("E51").Value = '  +1.18%  '
